$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain value updates (safe to let Excel infer the cell type; the
# resulting displayed value is identical to the source text).
$changes = @(
    @{ Cell = 'D2'; Value = '34.099.23' }
    @{ Cell = 'E2'; Value = '  +0.03%  ' }
    @{ Cell = 'D3'; Value = '1.790.53' }
    @{ Cell = 'E4'; Value = '  +0.06%  ' }
    @{ Cell = 'D5'; Value = '227.23' }
    @{ Cell = 'E5'; Value = '  +0.56%  ' }
    @{ Cell = 'D6'; Value = '0.547' }
    @{ Cell = 'E8'; Value = '  -1.85%  ' }
    @{ Cell = 'E9'; Value = '  +2.92%  ' }
    @{ Cell = 'D10'; Value = '0.0692' }
    @{ Cell = 'E10'; Value = '  -2.77%  ' }
    @{ Cell = 'E11'; Value = '  +0.36%  ' }
    @{ Cell = 'D12'; Value = '2.048.00' }
    @{ Cell = 'E13'; Value = '  +3.58%  ' }
    @{ Cell = 'D14'; Value = '1.792.04' }
    @{ Cell = 'E14'; Value = '  +0.40%  ' }
    @{ Cell = 'B15'; Value = 'Polygon' }
    @{ Cell = 'C15'; Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic' }
    @{ Cell = 'D15'; Value = '0.622' }
    @{ Cell = 'E15'; Value = '  +0.32%  ' }
    @{ Cell = 'B16'; Value = 'WrappedBTC' }
    @{ Cell = 'C16'; Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc' }
    @{ Cell = 'D16'; Value = '34.085.39' }
    @{ Cell = 'E16'; Value = '  +0.16%  ' }
    @{ Cell = 'D17'; Value = '4.19' }
    @{ Cell = 'E17'; Value = '  +0.75%  ' }
    @{ Cell = 'D18'; Value = '67.87' }
    @{ Cell = 'E18'; Value = '  -0.07%  ' }
    @{ Cell = 'E19'; Value = '  +0.36%  ' }
    @{ Cell = 'D20'; Value = '0.0₃0781' }
    @{ Cell = 'E20'; Value = '  -0.90%  ' }
    @{ Cell = 'D21'; Value = '10.93' }
    @{ Cell = 'E21'; Value = '  +1.52%  ' }
    @{ Cell = 'E22'; Value = '  +0.10%  ' }
    @{ Cell = 'E23'; Value = '  +0.68%  ' }
    @{ Cell = 'E24'; Value = '  -2.31%  ' }
    @{ Cell = 'E25'; Value = '  +0.99%  ' }
    @{ Cell = 'E26'; Value = '  +1.65%  ' }
    @{ Cell = 'D27'; Value = '16.29' }
    @{ Cell = 'E27'; Value = '  -0.24%  ' }
    @{ Cell = 'D28'; Value = '0.114' }
    @{ Cell = 'E28'; Value = '  +0.96%  ' }
    @{ Cell = 'E29'; Value = '  +0.22%  ' }
    @{ Cell = 'E30'; Value = '  +1.13%  ' }
    @{ Cell = 'E31'; Value = '  +1.46%  ' }
    @{ Cell = 'D32'; Value = '3.67' }
    @{ Cell = 'E32'; Value = '  +0.85%  ' }
    @{ Cell = 'E33'; Value = '  +3.07%  ' }
    @{ Cell = 'E34'; Value = '  +0.75%  ' }
    @{ Cell = 'D35'; Value = '1.445.13' }
    @{ Cell = 'E35'; Value = '  +3.93%  ' }
    @{ Cell = 'D36'; Value = '0.645' }
    @{ Cell = 'E36'; Value = '  -0.27%  ' }
    @{ Cell = 'E37'; Value = '  +2.47%  ' }
    @{ Cell = 'D38'; Value = '2.36' }
    @{ Cell = 'E38'; Value = '  +7.12%  ' }
    @{ Cell = 'E39'; Value = '  -1.43%  ' }
    @{ Cell = 'D40'; Value = '80.51' }
    @{ Cell = 'E40'; Value = '  +2.98%  ' }
    @{ Cell = 'E41'; Value = '  +0.42%  ' }
    @{ Cell = 'D42'; Value = '0.924' }
    @{ Cell = 'E42'; Value = '  +0.97%  ' }
    @{ Cell = 'E43'; Value = '  +0.10%  ' }
    @{ Cell = 'E44'; Value = '  +6.38%  ' }
    @{ Cell = 'B45'; Value = 'FraxShare' }
    @{ Cell = 'C45'; Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs' }
    @{ Cell = 'D45'; Value = '6.07' }
    @{ Cell = 'E45'; Value = '  +4.11%  ' }
    @{ Cell = 'B46'; Value = 'BabyDogeCoin' }
    @{ Cell = 'C46'; Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge' }
    @{ Cell = 'D46'; Value = '0.0₆0138' }
    @{ Cell = 'E46'; Value = '  -2.66%  ' }
    @{ Cell = 'B47'; Value = 'Kaspa' }
    @{ Cell = 'C47'; Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas' }
    @{ Cell = 'D47'; Value = '0.0509' }
    @{ Cell = 'E47'; Value = '  +2.18%  ' }
    @{ Cell = 'E48'; Value = '  -0.62%  ' }
    @{ Cell = 'D49'; Value = '107.73' }
    @{ Cell = 'E49'; Value = '  -0.42%  ' }
    @{ Cell = 'D50'; Value = '1.949.43' }
    @{ Cell = 'E50'; Value = '  +0.31%  ' }
    @{ Cell = 'E51'; Value = '  +0.11%  ' }
)

foreach ($change in $changes) {
    $ws.Range($change.Cell).Value = $change.Value
}

# A handful of price cells end in a significant trailing zero (e.g. "32.20").
# Assigning them through .Value would make Excel parse the text as a number
# and silently drop the trailing zero (32.20 -> 32.2), which would not match
# the scraped source text, so force those specific cells to keep their text
# formatting before writing the value.
$textProtectedChanges = @(
    @{ Cell = 'D8'; Value = '32.20' }
    @{ Cell = 'D19'; Value = '245.70' }
    @{ Cell = 'D25'; Value = '161.90' }
)

foreach ($change in $textProtectedChanges) {
    $cell = $ws.Range($change.Cell)
    $cell.NumberFormat = "@"
    $cell.Value = $change.Value
}

